$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I41").Value = 0
$ws.Range("L41").Value = 224.4
$ws.Range("K41").Value = 0
$ws.Range("J41").Value = 224.4
$ws.Range("N41").Value = -1104.4
$ws.Range("H41").Value = 224.4
$ws.Range("M41").ClearContents()
$ws.Range("K74").Value = 3600
$ws.Range("H74").Value = 4310.5713
$ws.Range("I74").Value = 3600
$ws.Range("J74").Value = 5258
$ws.Range("N74").Value = -7130
$ws.Range("M74").Value = -2664
$ws.Range("L74").Value = 5258
$ws.Range("M77").Value = -13320
$ws.Range("L77").Value = 26290
$ws.Range("J77").Value = 5258
$ws.Range("H77").Value = 4310.5713
$ws.Range("N77").Value = -35650
$ws.Range("K77").Value = 18000
$ws.Range("I77").Value = 3600
$ws.Range("N98").Value = -5996
$ws.Range("M98").Value = 1064.22223
$ws.Range("L98").Value = 3000
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 433.77777
$ws.Range("H98").Value = 690.4
$ws.Range("I98").Value = 433.77777
$ws.Range("K122").Value = 1301.33331
$ws.Range("M122").Value = 1148.66669
$ws.Range("N122").Value = -13900
$ws.Range("L122").Value = 9000
$ws.Range("J122").Value = 3000
$ws.Range("I122").Value = 433.77777
$ws.Range("H122").Value = 690.4
$ws.Range("H137").Value = 2838.9512
$ws.Range("I137").Value = 2096.2727
$ws.Range("K137").Value = 6288.8181
$ws.Range("M137").Value = -3738.8181
$ws.Range("K138").Value = 4193.3514
$ws.Range("J138").Value = 2605.5264
$ws.Range("H138").Value = 1807.5536
$ws.Range("N138").Value = -18096.5792
$ws.Range("L138").Value = 7816.5792
$ws.Range("I138").Value = 1397.7838
$ws.Range("M138").Value = 946.6486000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N2").Value = -1583
$ws.Range("H2").Value = 2164.7334
$ws.Range("I2").Value = 2289
$ws.Range("J2").Value = 1357
$ws.Range("K2").Value = 2289
$ws.Range("L2").Value = 1357
$ws.Range("M2").Value = -2176
$ws.Range("J61").Value = 3663.4546
$ws.Range("H61").Value = 2247
$ws.Range("L61").Value = 3663.4546
$ws.Range("N61").Value = -4087.4546
$ws.Range("H102").Value = 2136.375
$ws.Range("I102").Value = 2089.1904
$ws.Range("K102").Value = 2089.1904
$ws.Range("M102").Value = -467.1904
$ws.Range("K116").Value = 2289
$ws.Range("H116").Value = 2164.7334
$ws.Range("I116").Value = 2289
$ws.Range("L116").Value = 1357
$ws.Range("M116").Value = 5
$ws.Range("J116").Value = 1357
$ws.Range("N116").Value = -5945
$ws.Range("L132").Value = 17817.375
$ws.Range("H132").Value = 3278.8032
$ws.Range("K132").Value = 6998.7333
$ws.Range("J132").Value = 5939.125
$ws.Range("N132").Value = -22877.375
$ws.Range("M132").Value = -4468.7333
$ws.Range("I132").Value = 2332.9111
$ws.Range("L136").Value = 10990.3638
$ws.Range("J136").Value = 3663.4546
$ws.Range("H136").Value = 2247
$ws.Range("N136").Value = -16090.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N3").Value = -1585
$ws.Range("H3").Value = 2164.7334
$ws.Range("J3").Value = 1357
$ws.Range("K3").Value = 2289
$ws.Range("L3").Value = 1357
$ws.Range("I3").Value = 2289
$ws.Range("M3").Value = -2175
$ws.Range("H86").Value = 90911150
$ws.Range("K86").Value = 111113130
$ws.Range("I86").Value = 111113130
$ws.Range("M86").Value = -111112007
$ws.Range("H89").Value = 90911150
$ws.Range("K89").Value = 555565650
$ws.Range("M89").Value = -555560034
$ws.Range("I89").Value = 111113130
$ws.Range("J103").Value = 5000
$ws.Range("H103").Value = 5000
$ws.Range("L103").Value = 5000
$ws.Range("N103").Value = -7344
$ws.Range("J134").Value = 3921.3
$ws.Range("K134").Value = 6809.768999999999
$ws.Range("H134").Value = 2728.639
$ws.Range("I134").Value = 2269.923
$ws.Range("L134").Value = 11763.9
$ws.Range("N134").Value = -16833.9
$ws.Range("M134").Value = -4274.768999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7335.614
$ws.Range("M31").Value = -997.0416
$ws.Range("L31").Value = 14587.9
$ws.Range("J31").Value = 14587.9
$ws.Range("N31").Value = -15177.9
$ws.Range("K31").Value = 1292.0416
$ws.Range("I31").Value = 1292.0416
$ws.Range("L34").Value = 14587.9
$ws.Range("I34").Value = 1292.0416
$ws.Range("K34").Value = 1292.0416
$ws.Range("H34").Value = 7335.614
$ws.Range("M34").Value = -1090.0416
$ws.Range("N34").Value = -14991.9
$ws.Range("J34").Value = 14587.9
$ws.Range("L58").Value = 2757
$ws.Range("N58").Value = -3163
$ws.Range("K58").Value = 1454.8
$ws.Range("H58").Value = 1671.8334
$ws.Range("M58").Value = -1251.8
$ws.Range("I58").Value = 1454.8
$ws.Range("J58").Value = 2757
$ws.Range("K122").Value = 2615.4999
$ws.Range("M122").Value = -165.4998999999998
$ws.Range("N122").Value = -10547.0587
$ws.Range("L122").Value = 5647.0587
$ws.Range("J122").Value = 1882.3529
$ws.Range("I122").Value = 871.8333
$ws.Range("H122").Value = 1362.6571
$ws.Range("L132").Value = 71436714
$ws.Range("H132").Value = 6412230.5
$ws.Range("K132").Value = 5102.6844
$ws.Range("J132").Value = 23812238
$ws.Range("N132").Value = -71441774
$ws.Range("M132").Value = -2572.6844
$ws.Range("I132").Value = 1700.8948
$ws.Range("J134").Value = 7001.25
$ws.Range("K134").Value = 13000.0005
$ws.Range("H134").Value = 5857.857
$ws.Range("I134").Value = 4333.3335
$ws.Range("L134").Value = 21003.75
$ws.Range("N134").Value = -26073.75
$ws.Range("M134").Value = -10465.0005
$ws.Range("L136").Value = 8271
$ws.Range("J136").Value = 2757
$ws.Range("K136").Value = 4364.4
$ws.Range("H136").Value = 1671.8334
$ws.Range("I136").Value = 1454.8
$ws.Range("N136").Value = -13371
$ws.Range("M136").Value = -1814.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K113").Value = 1869.1875
$ws.Range("L113").Value = 4528.2498
$ws.Range("H113").Value = 1002.9286
$ws.Range("J113").Value = 1509.4166
$ws.Range("N113").Value = -8868.2498
$ws.Range("M113").Value = 300.8125
$ws.Range("I113").Value = 623.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L132").Value = 9128.25
$ws.Range("H132").Value = 2718.6785
$ws.Range("K132").Value = 7426.875
$ws.Range("J132").Value = 3042.75
$ws.Range("N132").Value = -14188.25
$ws.Range("M132").Value = -4896.875
$ws.Range("I132").Value = 2475.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J55").Value = 827.1429000000001
$ws.Range("L55").Value = 827.1429000000001
$ws.Range("M55").Value = -250.25
$ws.Range("N55").Value = -1173.1429
$ws.Range("I55").Value = 423.25
$ws.Range("H55").Value = 680.2727
$ws.Range("K55").Value = 423.25
$ws.Range("K132").Value = 5915.7999
$ws.Range("H132").Value = 2643
$ws.Range("M132").Value = -3385.7999
$ws.Range("I132").Value = 1971.9333
$ws.Range("L136").Value = 250015020
$ws.Range("J136").Value = 83338340
$ws.Range("K136").Value = 5559.706200000001
$ws.Range("H136").Value = 8774115
$ws.Range("I136").Value = 1853.2354
$ws.Range("N136").Value = -250020120
$ws.Range("M136").Value = -3009.706200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K132").Value = 4305.450000000001
$ws.Range("H132").Value = 6668455
$ws.Range("M132").Value = -1775.450000000001
$ws.Range("I132").Value = 1435.15
$ws.Range("L136").Value = 8994.643199999999
$ws.Range("J136").Value = 2998.2144
$ws.Range("K136").Value = 6197.6124
$ws.Range("H136").Value = 2355.9333
$ws.Range("I136").Value = 2065.8708
$ws.Range("N136").Value = -14094.6432
$ws.Range("M136").Value = -3647.6124
